# Update the 25 multiplication-fact answers in the table to the new
# problem/answer pairs, leaving all other document content untouched.
$d = $word.ActiveDocument

$d.Content.Find.Execute("14×81=1134", $true, $false, $false, $false, $false, $true, 1, $false, "79×94=7426", 2) | Out-Null
$d.Content.Find.Execute("92×92=8464", $true, $false, $false, $false, $false, $true, 1, $false, "57×62=3534", 2) | Out-Null
$d.Content.Find.Execute("64×96=6144", $true, $false, $false, $false, $false, $true, 1, $false, "58×92=5336", 2) | Out-Null
$d.Content.Find.Execute("37×13=481", $true, $false, $false, $false, $false, $true, 1, $false, "74×60=4440", 2) | Out-Null
$d.Content.Find.Execute("43×67=2881", $true, $false, $false, $false, $false, $true, 1, $false, "46×92=4232", 2) | Out-Null
$d.Content.Find.Execute("93×73=6789", $true, $false, $false, $false, $false, $true, 1, $false, "23×67=1541", 2) | Out-Null
$d.Content.Find.Execute("32×44=1408", $true, $false, $false, $false, $false, $true, 1, $false, "61×71=4331", 2) | Out-Null
$d.Content.Find.Execute("96×34=3264", $true, $false, $false, $false, $false, $true, 1, $false, "73×50=3650", 2) | Out-Null
$d.Content.Find.Execute("38×59=2242", $true, $false, $false, $false, $false, $true, 1, $false, "91×45=4095", 2) | Out-Null
$d.Content.Find.Execute("35×87=3045", $true, $false, $false, $false, $false, $true, 1, $false, "86×71=6106", 2) | Out-Null
$d.Content.Find.Execute("97×96=9312", $true, $false, $false, $false, $false, $true, 1, $false, "32×18=576", 2) | Out-Null
$d.Content.Find.Execute("81×48=3888", $true, $false, $false, $false, $false, $true, 1, $false, "12×48=576", 2) | Out-Null
$d.Content.Find.Execute("63×63=3969", $true, $false, $false, $false, $false, $true, 1, $false, "52×48=2496", 2) | Out-Null
$d.Content.Find.Execute("41×47=1927", $true, $false, $false, $false, $false, $true, 1, $false, "82×39=3198", 2) | Out-Null
$d.Content.Find.Execute("38×71=2698", $true, $false, $false, $false, $false, $true, 1, $false, "78×33=2574", 2) | Out-Null
$d.Content.Find.Execute("50×65=3250", $true, $false, $false, $false, $false, $true, 1, $false, "29×74=2146", 2) | Out-Null
$d.Content.Find.Execute("89×41=3649", $true, $false, $false, $false, $false, $true, 1, $false, "99×19=1881", 2) | Out-Null
$d.Content.Find.Execute("70×77=5390", $true, $false, $false, $false, $false, $true, 1, $false, "99×75=7425", 2) | Out-Null
$d.Content.Find.Execute("13×54=702", $true, $false, $false, $false, $false, $true, 1, $false, "20×30=600", 2) | Out-Null
$d.Content.Find.Execute("27×24=648", $true, $false, $false, $false, $false, $true, 1, $false, "84×25=2100", 2) | Out-Null
$d.Content.Find.Execute("56×55=3080", $true, $false, $false, $false, $false, $true, 1, $false, "39×36=1404", 2) | Out-Null
$d.Content.Find.Execute("78×95=7410", $true, $false, $false, $false, $false, $true, 1, $false, "82×56=4592", 2) | Out-Null
$d.Content.Find.Execute("38×53=2014", $true, $false, $false, $false, $false, $true, 1, $false, "43×48=2064", 2) | Out-Null
$d.Content.Find.Execute("69×23=1587", $true, $false, $false, $false, $false, $true, 1, $false, "56×72=4032", 2) | Out-Null
$d.Content.Find.Execute("22×15=330", $true, $false, $false, $false, $false, $true, 1, $false, "71×22=1562", 2) | Out-Null
